# Insert a new record row into the GESTION_TELECENTRO sheet.
# This pushes the former row 13 (and everything below it) down by one row,
# and fills the freshly inserted row 13 with the new case data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 13, shifting rows 13:50 -> 14:51
$ws.Rows.Item(13).Insert()

# Make sure the text-like columns (A:L) keep their original "text" storage
# (the sheet stores Caso/OT/etc. as text even though they look numeric),
# then write the values, then restore the default "Normal" style so no new
# cell formatting is introduced.
$textRange = $ws.Range("A13:L13")
$textRange.NumberFormat = "@"

$ws.Cells.Item(13, 1).Value  = "-148"
$ws.Cells.Item(13, 2).Value  = "10/21/2024"
$ws.Cells.Item(13, 3).Value  = "CARRANZA ANGEL J /ALT/ 2252"
$ws.Cells.Item(13, 4).Value  = "106594 - PALERMO"
$ws.Cells.Item(13, 5).Value  = "798311488"
$ws.Cells.Item(13, 6).Value  = "GESTION TELECENTRO"
$ws.Cells.Item(13, 7).Value  = "Pendiente"
$ws.Cells.Item(13, 8).Value  = "columna teco nodo tlc"
$ws.Cells.Item(13, 9).Value  = "0"
$ws.Cells.Item(13, 10).Value = "Cambio"
$ws.Cells.Item(13, 11).Value = "Nodo TLC"
$ws.Cells.Item(13, 12).Value = "Pasante"

$textRange.Style = "Normal"

# Numeric coordinate columns
$ws.Cells.Item(13, 13).Value = -58.433532
$ws.Cells.Item(13, 14).Value = -34.578254
